# Update column F ("dSF") values on Sheet1 to match re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = 4
    7  = 0
    8  = -1
    9  = 6
    16 = 3
    22 = 1
    27 = -1
    29 = 2
    31 = -1
    40 = 0
    42 = -2
    43 = 1
    51 = 0
    64 = -3
    67 = -4
    74 = -5
    79 = -14
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
